# GUN 3.0 FFS EN URENREGISTRATIE - WOLLAH
# Update week-2 (row 19/20) attendance data: add remark for Sam/Rief/Michiel,
# and mark Marc as fully present (with a "Marc Ziek" remark) for week 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 19: extend the existing remark in J19 ---
$ws.Range("J19").Value = "Sam Ziek, Rief thuissituatie, Michiel laptop"

# --- Row 20: Marc was fully present Mon-Thu (C:F) and Friday (I) ---
# Copy the "fully present" formatting (blue fill) from the matching
# columns of an existing fully-present row before setting the values,
# so the cell styles match the C/D/E/F/I pattern used elsewhere (e.g. row 12).
$ws.Range("C12").Copy()
$ws.Range("C20").PasteSpecial(-4122)

$ws.Range("D19").Copy()
$ws.Range("D20").PasteSpecial(-4122)

$ws.Range("D19").Copy()
$ws.Range("E20").PasteSpecial(-4122)

$ws.Range("D19").Copy()
$ws.Range("F20").PasteSpecial(-4122)

$ws.Range("I12").Copy()
$ws.Range("I20").PasteSpecial(-4122)

$ws.Range("C20").Value = 4
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = 4
$ws.Range("F20").Value = 4
$ws.Range("I20").Value = 4

$ws.Range("J20").Value = "Marc Ziek"

# --- Row 24 "Totaal lesuren" (column B) is a manually maintained weekly
# total, not a formula, so bump it to match the corrected week total. ---
$ws.Range("B24").Value = 8

# --- Scroll the view over to the right a bit and move the selection,
# matching where the author was working. ---
$ws.Range("M7").Select()
